# "wire_transfers" was a misnomer (see commit message) -- rename it to
# something that actually reflects the sheet's purpose, and since it's the
# sheet being discussed/edited right now, make it the active tab too.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wire_transfers")
$ws.Name = "currency conversion to EUR"
$ws.Activate()
